# Update schedule data for Linea 141 (18/01/2026 scrape run)
# New scrape timestamp: 10:56:15 (previous: 10:49:38)

$wb = $excel.ActiveWorkbook

$newScrapeTime = "10:56:15"

$dataLP1912 = @(
    @("05:57:04","06:09","10_OLMOS",12,"LP1912"),
    @("05:57:04","06:16","215A_EL PATO",19,"LP1912"),
    @("05:57:04","06:30","23_HERNANDEZ",33,"LP1912"),
    @("05:57:04","06:34","11_ETCHEVERRY",37,"LP1912"),
    @("06:34:35","06:36","11_ETCHEVERRY",2,"LP1912"),
    @("05:57:04","06:39","17X38_ROMERO",42,"LP1912"),
    @("05:57:04","06:41","16_SANTA ANA",44,"LP1912"),
    @("06:16:41","06:56","215A_EL PATO",40,"LP1912"),
    @("05:57:04","06:57","215A_EL PATO",60,"LP1912"),
    @("05:57:04","06:59","225_GOMEZ",62,"LP1912"),
    @("06:16:41","07:15","215C_EL PATO",59,"LP1912"),
    @("05:57:04","07:16","215C_EL PATO",79,"LP1912"),
    @("05:57:04","07:19","14_ABASTO",82,"LP1912"),
    @("06:16:41","07:20","16_SANTA ANA",64,"LP1912"),
    @("06:16:41","07:21","23_HERNANDEZ",65,"LP1912"),
    @("05:57:04","07:21","16_SANTA ANA",84,"LP1912"),
    @("05:57:04","07:22","23_HERNANDEZ",85,"LP1912"),
    @("05:57:04","07:29","17X38_ROMERO",92,"LP1912"),
    @("05:57:04","07:35","10_OLMOS",98,"LP1912"),
    @("06:16:41","07:36","27_EL RETIRO",80,"LP1912"),
    @("05:57:04","07:37","27_EL RETIRO",100,"LP1912"),
    @("06:34:35","07:43","215A_EL PATO",69,"LP1912"),
    @("06:46:20","07:44","215A_EL PATO",58,"LP1912"),
    @("07:49:32","07:49","215A_EL PATO",0,"LP1912"),
    @("05:57:04","07:55","14_ABASTO",118,"LP1912"),
    @("07:56:02","07:59","14_ABASTO",3,"LP1912"),
    @("06:16:41","08:00","17_ROMERO",104,"LP1912"),
    @("06:16:41","08:01","16_SANTA ANA",105,"LP1912"),
    @("06:34:35","08:06","23_HERNANDEZ",92,"LP1912"),
    @("07:56:02","08:07","23_HERNANDEZ",11,"LP1912"),
    @("06:16:41","08:11","10_OLMOS",115,"LP1912"),
    @("06:16:41","08:13","15X38_ABASTO",117,"LP1912"),
    @("06:34:35","08:29","15_ABASTO",115,"LP1912"),
    @("06:34:35","08:29","11_ETCHEVERRY",115,"LP1912"),
    @("08:28:52","08:32","11_ETCHEVERRY",4,"LP1912"),
    @("08:38:24","08:40","10_OLMOS",2,"LP1912"),
    @("06:46:20","08:41","16_P MOR-SANTA ANA",115,"LP1912"),
    @("06:53:44","08:43","215C_EL PATO",110,"LP1912"),
    @("06:46:20","08:44","215C_EL PATO",118,"LP1912"),
    @("08:45:31","08:45","215C_EL PATO",0,"LP1912"),
    @("07:56:02","08:45","23_HERNANDEZ",49,"LP1912"),
    @("07:49:32","08:47","23_HERNANDEZ",58,"LP1912"),
    @("07:38:09","08:51","23_HERNANDEZ",73,"LP1912"),
    @("08:52:40","08:52","215B_EL PATO",0,"LP1912"),
    @("07:13:03","08:52","23_HERNANDEZ",99,"LP1912"),
    @("07:49:32","08:53","215B_EL PATO",64,"LP1912"),
    @("07:13:03","08:54","215B_EL PATO",101,"LP1912"),
    @("08:28:52","08:57","215A_EL PATO",29,"LP1912"),
    @("07:13:03","08:58","215A_EL PATO",105,"LP1912"),
    @("08:28:52","09:04","10_OLMOS",36,"LP1912"),
    @("08:11:18","09:05","10_OLMOS",54,"LP1912"),
    @("07:38:09","09:06","16_SANTA ANA",88,"LP1912"),
    @("07:56:02","09:11","27_EL RETIRO",75,"LP1912"),
    @("07:49:32","09:12","27_EL RETIRO",83,"LP1912"),
    @("07:38:09","09:14","27_EL RETIRO",96,"LP1912"),
    @("08:28:52","09:16","27_EL RETIRO",48,"LP1912"),
    @("08:38:24","09:17","27_EL RETIRO",39,"LP1912"),
    @("07:49:32","09:17","14_ABASTO",88,"LP1912"),
    @("07:38:09","09:18","15X38_ABASTO",100,"LP1912"),
    @("07:38:09","09:18","14_ABASTO",100,"LP1912"),
    @("08:28:52","09:28","10_OLMOS",60,"LP1912"),
    @("08:11:18","09:28","23_HERNANDEZ",77,"LP1912"),
    @("07:38:09","09:29","10_OLMOS",111,"LP1912"),
    @("08:11:18","09:31","16_SANTA ANA",80,"LP1912"),
    @("08:28:52","09:32","23_HERNANDEZ",64,"LP1912"),
    @("08:38:24","09:34","23_HERNANDEZ",56,"LP1912"),
    @("08:52:40","09:36","23_HERNANDEZ",44,"LP1912"),
    @("07:49:32","09:39","15_ABASTO",110,"LP1912"),
    @("07:49:32","09:41","11_ETCHEVERRY",112,"LP1912"),
    @("07:56:02","09:42","11_ETCHEVERRY",106,"LP1912"),
    @("07:49:32","09:43","16_P MOR-SANTA ANA",114,"LP1912"),
    @("08:11:18","09:53","10_OLMOS",102,"LP1912"),
    @("08:28:52","09:58","215C_EL PATO",90,"LP1912"),
    @("08:11:18","09:59","215C_EL PATO",108,"LP1912"),
    @("08:28:52","10:05","14_ABASTO",97,"LP1912"),
    @("08:11:18","10:06","14_ABASTO",115,"LP1912"),
    @("08:28:52","10:13","17X38_ROMERO",105,"LP1912"),
    @("09:22:34","10:21","23_HERNANDEZ",59,"LP1912"),
    @("10:04:30","10:22","23_HERNANDEZ",18,"LP1912"),
    @("09:22:34","10:25","16_SANTA ANA",63,"LP1912"),
    @("08:38:24","10:29","15_ABASTO",111,"LP1912"),
    @("10:04:30","10:29","14_ABASTO",25,"LP1912"),
    @("10:36:50","10:43","11X44_ETCHEVERRY",7,"LP1912"),
    @("08:45:31","10:44","11X44_ETCHEVERRY",119,"LP1912"),
    @("08:52:40","10:46","15_P INDUSTRIAL",114,"LP1912"),
    @("09:22:34","10:53","27_EL RETIRO",91,"LP1912"),
    @("10:36:50","10:55","16_SANTA ANA",19,"LP1912"),
    @("10:04:30","10:56","27_EL RETIRO",52,"LP1912"),
    @("09:22:34","10:57","10_OLMOS",95,"LP1912"),
    @("10:56:15","10:57","27_EL RETIRO",1,"LP1912"),
    @("10:04:30","10:59","10_OLMOS",55,"LP1912"),
    @("09:22:34","11:01","81_EL PELIGRO",99,"LP1912"),
    @("10:04:30","11:03","23_HERNANDEZ",59,"LP1912"),
    @("10:36:50","11:06","23_HERNANDEZ",30,"LP1912"),
    @("09:22:34","11:10","16_P MOR-SANTA ANA",108,"LP1912"),
    @("09:22:34","11:14","14_ABASTO",112,"LP1912"),
    @("09:22:34","11:15","15X38_ABASTO",113,"LP1912"),
    @("10:56:15","11:15","14_ABASTO",19,"LP1912"),
    @("10:36:50","11:25","16_SANTA ANA",49,"LP1912"),
    @("10:04:30","11:29","10_OLMOS",85,"LP1912"),
    @("10:04:30","11:29","16_SANTA ANA",85,"LP1912"),
    @("10:36:50","11:30","215C_EL PATO",54,"LP1912"),
    @("10:04:30","11:31","215C_EL PATO",87,"LP1912"),
    @("10:04:30","11:41","215B_EL PATO",97,"LP1912"),
    @("10:56:15","11:42","215B_EL PATO",46,"LP1912"),
    @("10:04:30","11:45","15X38_ABASTO",101,"LP1912"),
    @("10:56:15","11:46","23_HERNANDEZ",50,"LP1912"),
    @("10:49:38","11:47","23_HERNANDEZ",58,"LP1912"),
    @("10:36:50","11:48","23_HERNANDEZ",72,"LP1912"),
    @("10:36:50","11:52","225_GOMEZ",76,"LP1912"),
    @("10:04:30","11:53","225_GOMEZ",109,"LP1912"),
    @("10:04:30","11:58","17_ROMERO",114,"LP1912"),
    @("10:36:50","12:05","11_ETCHEVERRY",89,"LP1912"),
    @("10:56:15","12:06","11_ETCHEVERRY",70,"LP1912"),
    @("10:36:50","12:10","15_ABASTO",94,"LP1912"),
    @("10:36:50","12:10","16_P MOR-SANTA ANA",94,"LP1912"),
    @("10:36:50","12:21","215C_EL PATO",105,"LP1912"),
    @("10:56:15","12:22","215C_EL PATO",86,"LP1912"),
    @("10:36:50","12:32","14_ABASTO",116,"LP1912"),
    @("10:56:15","12:33","14_ABASTO",97,"LP1912"),
    @("10:56:15","12:33","27_EL RETIRO",97,"LP1912"),
    @("10:36:50","12:34","15_ABASTO",118,"LP1912"),
    @("10:49:38","12:36","27_EL RETIRO",107,"LP1912"),
    @("10:49:38","12:48","16_SANTA ANA",119,"LP1912")
)

$dataLP1912_215 = @(
    @("05:57:04","06:16","215A_EL PATO",19,"LP1912"),
    @("06:16:41","06:56","215A_EL PATO",40,"LP1912"),
    @("05:57:04","06:57","215A_EL PATO",60,"LP1912"),
    @("06:16:41","07:15","215C_EL PATO",59,"LP1912"),
    @("05:57:04","07:16","215C_EL PATO",79,"LP1912"),
    @("06:34:35","07:43","215A_EL PATO",69,"LP1912"),
    @("06:46:20","07:44","215A_EL PATO",58,"LP1912"),
    @("07:49:32","07:49","215A_EL PATO",0,"LP1912"),
    @("06:53:44","08:43","215C_EL PATO",110,"LP1912"),
    @("06:46:20","08:44","215C_EL PATO",118,"LP1912"),
    @("08:45:31","08:45","215C_EL PATO",0,"LP1912"),
    @("08:52:40","08:52","215B_EL PATO",0,"LP1912"),
    @("07:49:32","08:53","215B_EL PATO",64,"LP1912"),
    @("07:13:03","08:54","215B_EL PATO",101,"LP1912"),
    @("08:28:52","08:57","215A_EL PATO",29,"LP1912"),
    @("07:13:03","08:58","215A_EL PATO",105,"LP1912"),
    @("08:28:52","09:58","215C_EL PATO",90,"LP1912"),
    @("08:11:18","09:59","215C_EL PATO",108,"LP1912"),
    @("10:36:50","11:30","215C_EL PATO",54,"LP1912"),
    @("10:04:30","11:31","215C_EL PATO",87,"LP1912"),
    @("10:04:30","11:41","215B_EL PATO",97,"LP1912"),
    @("10:56:15","11:42","215B_EL PATO",46,"LP1912"),
    @("10:36:50","12:21","215C_EL PATO",105,"LP1912"),
    @("10:56:15","12:22","215C_EL PATO",86,"LP1912")
)


# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: " + $newScrapeTime
$ws1.Range("A3").Value = "Total filas: " + $dataLP1912.Count

for ($i = 0; $i -lt $dataLP1912.Count; $i++) {
    $row = 6 + $i
    $rec = $dataLP1912[$i]
    $ws1.Cells.Item($row, 1).Value = $rec[0]
    $ws1.Cells.Item($row, 2).Value = $rec[1]
    $ws1.Cells.Item($row, 3).Value = $rec[2]
    $ws1.Cells.Item($row, 4).Value = $rec[3]
    $ws1.Cells.Item($row, 5).Value = $rec[4]
}

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: " + $newScrapeTime
$ws2.Range("A3").Value = "Total filas: " + $dataLP1912_215.Count

for ($i = 0; $i -lt $dataLP1912_215.Count; $i++) {
    $row = 6 + $i
    $rec = $dataLP1912_215[$i]
    $ws2.Cells.Item($row, 1).Value = $rec[0]
    $ws2.Cells.Item($row, 2).Value = $rec[1]
    $ws2.Cells.Item($row, 3).Value = $rec[2]
    $ws2.Cells.Item($row, 4).Value = $rec[3]
    $ws2.Cells.Item($row, 5).Value = $rec[4]
}

# ---- Sheet 3: 6203-6173 ----
# Only the "last updated" timestamp changes here; the data rows are unchanged.
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: " + $newScrapeTime
